$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.699.66"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "2.471.84"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.89"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0852"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").Value = "2.853.21"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.47"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "2.472.03"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("D18").Value = "41.630.05"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").Value = "0.0₃0939"
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.43"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("E22").Value = "  -4.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.03"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.93"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.81"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.72"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.37"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.30"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.42"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.65%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0762"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.55"
$ws.Range("D35").ClearFormats()
$ws.Range("E36").Value = "  -2.96%  "
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.89"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.83"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.98"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("E42").Value = "  -3.84%  "
$ws.Range("D43").Value = "2.000.84"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.67"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("E46").Value = "  -0.11%  "
$ws.Range("E47").Value = "  +4.08%  "
$ws.Range("D48").Value = "2.733.56"
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.25"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.71"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.39%  "
$ws.Range("E51").Value = "  -0.71%  "
